# Edit script: applies the "Add files via upload" changes to
# dati/provincia/Vibo Valentia/Vibo Valentia.xlsx
#
# 1) Two data corrections in "Nuovi casi" (sheet1): C426 22->17, C494 0->-2
# 2) New daily rows 509-521 (C=new case counts, D=7-day rolling average)
#    added to every sheet (Nuovi casi, Deceduti, Dimessi Guariti, Ricoveri).
# 3) New trailing date-only rows 522-541 (column A only) on every sheet.
# 4) Selection/scroll-position bookkeeping to match the saved view state.

$wb = $excel.ActiveWorkbook
$dateFmt = "dd/mm/yyyy"
$numFmt = "#,##0"

# ---------------------------------------------------------------
# 1) Data corrections on "Nuovi casi"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(426,3).Value = 17
$ws1.Cells.Item(494,3).Value = -2

# ---------------------------------------------------------------
# 2) New daily rows 509-521 on sheet1
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(509,3).Value = 0
$ws.Cells.Item(509,4).Formula = "=AVERAGE(C503:C509)"
$ws.Cells.Item(509,4).Font.Color = 0
$ws.Cells.Item(509,4).NumberFormat = $numFmt
$ws.Cells.Item(510,3).Value = 8
$ws.Cells.Item(510,4).Formula = "=AVERAGE(C504:C510)"
$ws.Cells.Item(510,4).Font.Color = 0
$ws.Cells.Item(510,4).NumberFormat = $numFmt
$ws.Cells.Item(511,3).Value = 8
$ws.Cells.Item(511,4).Formula = "=AVERAGE(C505:C511)"
$ws.Cells.Item(511,4).Font.Color = 0
$ws.Cells.Item(511,4).NumberFormat = $numFmt
$ws.Cells.Item(512,3).Value = 0
$ws.Cells.Item(512,4).Formula = "=AVERAGE(C506:C512)"
$ws.Cells.Item(512,4).Font.Color = 0
$ws.Cells.Item(512,4).NumberFormat = $numFmt
$ws.Cells.Item(513,3).Value = 8
$ws.Cells.Item(513,4).Formula = "=AVERAGE(C507:C513)"
$ws.Cells.Item(513,4).Font.Color = 0
$ws.Cells.Item(513,4).NumberFormat = $numFmt
$ws.Cells.Item(514,3).Value = 3
$ws.Cells.Item(514,4).Formula = "=AVERAGE(C508:C514)"
$ws.Cells.Item(514,4).Font.Color = 0
$ws.Cells.Item(514,4).NumberFormat = $numFmt
$ws.Cells.Item(515,3).Value = 9
$ws.Cells.Item(515,4).Formula = "=AVERAGE(C509:C515)"
$ws.Cells.Item(515,4).Font.Color = 0
$ws.Cells.Item(515,4).NumberFormat = $numFmt
$ws.Cells.Item(516,3).Value = 11
$ws.Cells.Item(516,4).Formula = "=AVERAGE(C510:C516)"
$ws.Cells.Item(516,4).Font.Color = 0
$ws.Cells.Item(516,4).NumberFormat = $numFmt
$ws.Cells.Item(517,3).Value = 23
$ws.Cells.Item(517,4).Formula = "=AVERAGE(C511:C517)"
$ws.Cells.Item(517,4).Font.Color = 0
$ws.Cells.Item(517,4).NumberFormat = $numFmt
$ws.Cells.Item(518,3).Value = 24
$ws.Cells.Item(518,4).Formula = "=AVERAGE(C512:C518)"
$ws.Cells.Item(518,4).Font.Color = 0
$ws.Cells.Item(518,4).NumberFormat = $numFmt
$ws.Cells.Item(519,3).Value = 0
$ws.Cells.Item(519,4).Formula = "=AVERAGE(C513:C519)"
$ws.Cells.Item(519,4).Font.Color = 0
$ws.Cells.Item(519,4).NumberFormat = $numFmt
$ws.Cells.Item(520,3).Value = 18
$ws.Cells.Item(520,4).Formula = "=AVERAGE(C514:C520)"
$ws.Cells.Item(520,4).Font.Color = 0
$ws.Cells.Item(520,4).NumberFormat = $numFmt
$ws.Cells.Item(521,3).Value = 15
$ws.Cells.Item(521,4).Formula = "=AVERAGE(C515:C521)"
$ws.Cells.Item(521,4).Font.Color = 0
$ws.Cells.Item(521,4).NumberFormat = $numFmt

# ---------------------------------------------------------------
# 3) Trailing date-only rows 522-541 on sheet1
# ---------------------------------------------------------------
$ws.Cells.Item(522,1).Value = 44420
$ws.Cells.Item(522,1).NumberFormat = $dateFmt
$ws.Cells.Item(523,1).Value = 44421
$ws.Cells.Item(523,1).NumberFormat = $dateFmt
$ws.Cells.Item(524,1).Value = 44422
$ws.Cells.Item(524,1).NumberFormat = $dateFmt
$ws.Cells.Item(525,1).Value = 44423
$ws.Cells.Item(525,1).NumberFormat = $dateFmt
$ws.Cells.Item(526,1).Value = 44424
$ws.Cells.Item(526,1).NumberFormat = $dateFmt
$ws.Cells.Item(527,1).Value = 44425
$ws.Cells.Item(527,1).NumberFormat = $dateFmt
$ws.Cells.Item(528,1).Value = 44426
$ws.Cells.Item(528,1).NumberFormat = $dateFmt
$ws.Cells.Item(529,1).Value = 44427
$ws.Cells.Item(529,1).NumberFormat = $dateFmt
$ws.Cells.Item(530,1).Value = 44428
$ws.Cells.Item(530,1).NumberFormat = $dateFmt
$ws.Cells.Item(531,1).Value = 44429
$ws.Cells.Item(531,1).NumberFormat = $dateFmt
$ws.Cells.Item(532,1).Value = 44430
$ws.Cells.Item(532,1).NumberFormat = $dateFmt
$ws.Cells.Item(533,1).Value = 44431
$ws.Cells.Item(533,1).NumberFormat = $dateFmt
$ws.Cells.Item(534,1).Value = 44432
$ws.Cells.Item(534,1).NumberFormat = $dateFmt
$ws.Cells.Item(535,1).Value = 44433
$ws.Cells.Item(535,1).NumberFormat = $dateFmt
$ws.Cells.Item(536,1).Value = 44434
$ws.Cells.Item(536,1).NumberFormat = $dateFmt
$ws.Cells.Item(537,1).Value = 44435
$ws.Cells.Item(537,1).NumberFormat = $dateFmt
$ws.Cells.Item(538,1).Value = 44436
$ws.Cells.Item(538,1).NumberFormat = $dateFmt
$ws.Cells.Item(539,1).Value = 44437
$ws.Cells.Item(539,1).NumberFormat = $dateFmt
$ws.Cells.Item(540,1).Value = 44438
$ws.Cells.Item(540,1).NumberFormat = $dateFmt
$ws.Cells.Item(541,1).Value = 44439
$ws.Cells.Item(541,1).NumberFormat = $dateFmt

# ---------------------------------------------------------------
# 2) New daily rows 509-521 on sheet2
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(509,3).Value = 0
$ws.Cells.Item(509,4).Formula = "=AVERAGE(C503:C509)"
$ws.Cells.Item(509,4).Font.Color = 0
$ws.Cells.Item(509,4).NumberFormat = $numFmt
$ws.Cells.Item(510,3).Value = 0
$ws.Cells.Item(510,4).Formula = "=AVERAGE(C504:C510)"
$ws.Cells.Item(510,4).Font.Color = 0
$ws.Cells.Item(510,4).NumberFormat = $numFmt
$ws.Cells.Item(511,3).Value = 0
$ws.Cells.Item(511,4).Formula = "=AVERAGE(C505:C511)"
$ws.Cells.Item(511,4).Font.Color = 0
$ws.Cells.Item(511,4).NumberFormat = $numFmt
$ws.Cells.Item(512,3).Value = 0
$ws.Cells.Item(512,4).Formula = "=AVERAGE(C506:C512)"
$ws.Cells.Item(512,4).Font.Color = 0
$ws.Cells.Item(512,4).NumberFormat = $numFmt
$ws.Cells.Item(513,3).Value = 0
$ws.Cells.Item(513,4).Formula = "=AVERAGE(C507:C513)"
$ws.Cells.Item(513,4).Font.Color = 0
$ws.Cells.Item(513,4).NumberFormat = $numFmt
$ws.Cells.Item(514,3).Value = 0
$ws.Cells.Item(514,4).Formula = "=AVERAGE(C508:C514)"
$ws.Cells.Item(514,4).Font.Color = 0
$ws.Cells.Item(514,4).NumberFormat = $numFmt
$ws.Cells.Item(515,3).Value = 0
$ws.Cells.Item(515,4).Formula = "=AVERAGE(C509:C515)"
$ws.Cells.Item(515,4).Font.Color = 0
$ws.Cells.Item(515,4).NumberFormat = $numFmt
$ws.Cells.Item(516,3).Value = 0
$ws.Cells.Item(516,4).Formula = "=AVERAGE(C510:C516)"
$ws.Cells.Item(516,4).Font.Color = 0
$ws.Cells.Item(516,4).NumberFormat = $numFmt
$ws.Cells.Item(517,3).Value = 0
$ws.Cells.Item(517,4).Formula = "=AVERAGE(C511:C517)"
$ws.Cells.Item(517,4).Font.Color = 0
$ws.Cells.Item(517,4).NumberFormat = $numFmt
$ws.Cells.Item(518,3).Value = 0
$ws.Cells.Item(518,4).Formula = "=AVERAGE(C512:C518)"
$ws.Cells.Item(518,4).Font.Color = 0
$ws.Cells.Item(518,4).NumberFormat = $numFmt
$ws.Cells.Item(519,3).Value = 0
$ws.Cells.Item(519,4).Formula = "=AVERAGE(C513:C519)"
$ws.Cells.Item(519,4).Font.Color = 0
$ws.Cells.Item(519,4).NumberFormat = $numFmt
$ws.Cells.Item(520,3).Value = 0
$ws.Cells.Item(520,4).Formula = "=AVERAGE(C514:C520)"
$ws.Cells.Item(520,4).Font.Color = 0
$ws.Cells.Item(520,4).NumberFormat = $numFmt
$ws.Cells.Item(521,3).Value = 0
$ws.Cells.Item(521,4).Formula = "=AVERAGE(C515:C521)"
$ws.Cells.Item(521,4).Font.Color = 0
$ws.Cells.Item(521,4).NumberFormat = $numFmt

# ---------------------------------------------------------------
# 3) Trailing date-only rows 522-541 on sheet2
# ---------------------------------------------------------------
$ws.Cells.Item(522,1).Value = 44420
$ws.Cells.Item(522,1).NumberFormat = $dateFmt
$ws.Cells.Item(523,1).Value = 44421
$ws.Cells.Item(523,1).NumberFormat = $dateFmt
$ws.Cells.Item(524,1).Value = 44422
$ws.Cells.Item(524,1).NumberFormat = $dateFmt
$ws.Cells.Item(525,1).Value = 44423
$ws.Cells.Item(525,1).NumberFormat = $dateFmt
$ws.Cells.Item(526,1).Value = 44424
$ws.Cells.Item(526,1).NumberFormat = $dateFmt
$ws.Cells.Item(527,1).Value = 44425
$ws.Cells.Item(527,1).NumberFormat = $dateFmt
$ws.Cells.Item(528,1).Value = 44426
$ws.Cells.Item(528,1).NumberFormat = $dateFmt
$ws.Cells.Item(529,1).Value = 44427
$ws.Cells.Item(529,1).NumberFormat = $dateFmt
$ws.Cells.Item(530,1).Value = 44428
$ws.Cells.Item(530,1).NumberFormat = $dateFmt
$ws.Cells.Item(531,1).Value = 44429
$ws.Cells.Item(531,1).NumberFormat = $dateFmt
$ws.Cells.Item(532,1).Value = 44430
$ws.Cells.Item(532,1).NumberFormat = $dateFmt
$ws.Cells.Item(533,1).Value = 44431
$ws.Cells.Item(533,1).NumberFormat = $dateFmt
$ws.Cells.Item(534,1).Value = 44432
$ws.Cells.Item(534,1).NumberFormat = $dateFmt
$ws.Cells.Item(535,1).Value = 44433
$ws.Cells.Item(535,1).NumberFormat = $dateFmt
$ws.Cells.Item(536,1).Value = 44434
$ws.Cells.Item(536,1).NumberFormat = $dateFmt
$ws.Cells.Item(537,1).Value = 44435
$ws.Cells.Item(537,1).NumberFormat = $dateFmt
$ws.Cells.Item(538,1).Value = 44436
$ws.Cells.Item(538,1).NumberFormat = $dateFmt
$ws.Cells.Item(539,1).Value = 44437
$ws.Cells.Item(539,1).NumberFormat = $dateFmt
$ws.Cells.Item(540,1).Value = 44438
$ws.Cells.Item(540,1).NumberFormat = $dateFmt
$ws.Cells.Item(541,1).Value = 44439
$ws.Cells.Item(541,1).NumberFormat = $dateFmt

# ---------------------------------------------------------------
# 2) New daily rows 509-521 on sheet3
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(509,3).Value = 2
$ws.Cells.Item(509,4).Formula = "=AVERAGE(C503:C509)"
$ws.Cells.Item(509,4).Font.Color = 0
$ws.Cells.Item(509,4).NumberFormat = $numFmt
$ws.Cells.Item(510,3).Value = 2
$ws.Cells.Item(510,4).Formula = "=AVERAGE(C504:C510)"
$ws.Cells.Item(510,4).Font.Color = 0
$ws.Cells.Item(510,4).NumberFormat = $numFmt
$ws.Cells.Item(511,3).Value = 0
$ws.Cells.Item(511,4).Formula = "=AVERAGE(C505:C511)"
$ws.Cells.Item(511,4).Font.Color = 0
$ws.Cells.Item(511,4).NumberFormat = $numFmt
$ws.Cells.Item(512,3).Value = 1
$ws.Cells.Item(512,4).Formula = "=AVERAGE(C506:C512)"
$ws.Cells.Item(512,4).Font.Color = 0
$ws.Cells.Item(512,4).NumberFormat = $numFmt
$ws.Cells.Item(513,3).Value = 3
$ws.Cells.Item(513,4).Formula = "=AVERAGE(C507:C513)"
$ws.Cells.Item(513,4).Font.Color = 0
$ws.Cells.Item(513,4).NumberFormat = $numFmt
$ws.Cells.Item(514,3).Value = 4
$ws.Cells.Item(514,4).Formula = "=AVERAGE(C508:C514)"
$ws.Cells.Item(514,4).Font.Color = 0
$ws.Cells.Item(514,4).NumberFormat = $numFmt
$ws.Cells.Item(515,3).Value = 3
$ws.Cells.Item(515,4).Formula = "=AVERAGE(C509:C515)"
$ws.Cells.Item(515,4).Font.Color = 0
$ws.Cells.Item(515,4).NumberFormat = $numFmt
$ws.Cells.Item(516,3).Value = 3
$ws.Cells.Item(516,4).Formula = "=AVERAGE(C510:C516)"
$ws.Cells.Item(516,4).Font.Color = 0
$ws.Cells.Item(516,4).NumberFormat = $numFmt
$ws.Cells.Item(517,3).Value = 2
$ws.Cells.Item(517,4).Formula = "=AVERAGE(C511:C517)"
$ws.Cells.Item(517,4).Font.Color = 0
$ws.Cells.Item(517,4).NumberFormat = $numFmt
$ws.Cells.Item(518,3).Value = 6
$ws.Cells.Item(518,4).Formula = "=AVERAGE(C512:C518)"
$ws.Cells.Item(518,4).Font.Color = 0
$ws.Cells.Item(518,4).NumberFormat = $numFmt
$ws.Cells.Item(519,3).Value = 1
$ws.Cells.Item(519,4).Formula = "=AVERAGE(C513:C519)"
$ws.Cells.Item(519,4).Font.Color = 0
$ws.Cells.Item(519,4).NumberFormat = $numFmt
$ws.Cells.Item(520,3).Value = 0
$ws.Cells.Item(520,4).Formula = "=AVERAGE(C514:C520)"
$ws.Cells.Item(520,4).Font.Color = 0
$ws.Cells.Item(520,4).NumberFormat = $numFmt
$ws.Cells.Item(521,3).Value = 3
$ws.Cells.Item(521,4).Formula = "=AVERAGE(C515:C521)"
$ws.Cells.Item(521,4).Font.Color = 0
$ws.Cells.Item(521,4).NumberFormat = $numFmt

# ---------------------------------------------------------------
# 3) Trailing date-only rows 522-541 on sheet3
# ---------------------------------------------------------------
$ws.Cells.Item(522,1).Value = 44420
$ws.Cells.Item(522,1).NumberFormat = $dateFmt
$ws.Cells.Item(523,1).Value = 44421
$ws.Cells.Item(523,1).NumberFormat = $dateFmt
$ws.Cells.Item(524,1).Value = 44422
$ws.Cells.Item(524,1).NumberFormat = $dateFmt
$ws.Cells.Item(525,1).Value = 44423
$ws.Cells.Item(525,1).NumberFormat = $dateFmt
$ws.Cells.Item(526,1).Value = 44424
$ws.Cells.Item(526,1).NumberFormat = $dateFmt
$ws.Cells.Item(527,1).Value = 44425
$ws.Cells.Item(527,1).NumberFormat = $dateFmt
$ws.Cells.Item(528,1).Value = 44426
$ws.Cells.Item(528,1).NumberFormat = $dateFmt
$ws.Cells.Item(529,1).Value = 44427
$ws.Cells.Item(529,1).NumberFormat = $dateFmt
$ws.Cells.Item(530,1).Value = 44428
$ws.Cells.Item(530,1).NumberFormat = $dateFmt
$ws.Cells.Item(531,1).Value = 44429
$ws.Cells.Item(531,1).NumberFormat = $dateFmt
$ws.Cells.Item(532,1).Value = 44430
$ws.Cells.Item(532,1).NumberFormat = $dateFmt
$ws.Cells.Item(533,1).Value = 44431
$ws.Cells.Item(533,1).NumberFormat = $dateFmt
$ws.Cells.Item(534,1).Value = 44432
$ws.Cells.Item(534,1).NumberFormat = $dateFmt
$ws.Cells.Item(535,1).Value = 44433
$ws.Cells.Item(535,1).NumberFormat = $dateFmt
$ws.Cells.Item(536,1).Value = 44434
$ws.Cells.Item(536,1).NumberFormat = $dateFmt
$ws.Cells.Item(537,1).Value = 44435
$ws.Cells.Item(537,1).NumberFormat = $dateFmt
$ws.Cells.Item(538,1).Value = 44436
$ws.Cells.Item(538,1).NumberFormat = $dateFmt
$ws.Cells.Item(539,1).Value = 44437
$ws.Cells.Item(539,1).NumberFormat = $dateFmt
$ws.Cells.Item(540,1).Value = 44438
$ws.Cells.Item(540,1).NumberFormat = $dateFmt
$ws.Cells.Item(541,1).Value = 44439
$ws.Cells.Item(541,1).NumberFormat = $dateFmt

# ---------------------------------------------------------------
# 2) New daily rows 509-521 on sheet4
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(509,3).Value = 3
$ws.Cells.Item(509,4).Formula = "=AVERAGE(C503:C509)"
$ws.Cells.Item(509,4).Font.Color = 0
$ws.Cells.Item(509,4).NumberFormat = $numFmt
$ws.Cells.Item(510,3).Value = 3
$ws.Cells.Item(510,4).Formula = "=AVERAGE(C504:C510)"
$ws.Cells.Item(510,4).Font.Color = 0
$ws.Cells.Item(510,4).NumberFormat = $numFmt
$ws.Cells.Item(511,3).Value = 3
$ws.Cells.Item(511,4).Formula = "=AVERAGE(C505:C511)"
$ws.Cells.Item(511,4).Font.Color = 0
$ws.Cells.Item(511,4).NumberFormat = $numFmt
$ws.Cells.Item(512,3).Value = 2
$ws.Cells.Item(512,4).Formula = "=AVERAGE(C506:C512)"
$ws.Cells.Item(512,4).Font.Color = 0
$ws.Cells.Item(512,4).NumberFormat = $numFmt
$ws.Cells.Item(513,3).Value = 2
$ws.Cells.Item(513,4).Formula = "=AVERAGE(C507:C513)"
$ws.Cells.Item(513,4).Font.Color = 0
$ws.Cells.Item(513,4).NumberFormat = $numFmt
$ws.Cells.Item(514,3).Value = 1
$ws.Cells.Item(514,4).Formula = "=AVERAGE(C508:C514)"
$ws.Cells.Item(514,4).Font.Color = 0
$ws.Cells.Item(514,4).NumberFormat = $numFmt
$ws.Cells.Item(515,3).Value = 1
$ws.Cells.Item(515,4).Formula = "=AVERAGE(C509:C515)"
$ws.Cells.Item(515,4).Font.Color = 0
$ws.Cells.Item(515,4).NumberFormat = $numFmt
$ws.Cells.Item(516,3).Value = 2
$ws.Cells.Item(516,4).Formula = "=AVERAGE(C510:C516)"
$ws.Cells.Item(516,4).Font.Color = 0
$ws.Cells.Item(516,4).NumberFormat = $numFmt
$ws.Cells.Item(517,3).Value = 2
$ws.Cells.Item(517,4).Formula = "=AVERAGE(C511:C517)"
$ws.Cells.Item(517,4).Font.Color = 0
$ws.Cells.Item(517,4).NumberFormat = $numFmt
$ws.Cells.Item(518,3).Value = 4
$ws.Cells.Item(518,4).Formula = "=AVERAGE(C512:C518)"
$ws.Cells.Item(518,4).Font.Color = 0
$ws.Cells.Item(518,4).NumberFormat = $numFmt
$ws.Cells.Item(519,3).Value = 4
$ws.Cells.Item(519,4).Formula = "=AVERAGE(C513:C519)"
$ws.Cells.Item(519,4).Font.Color = 0
$ws.Cells.Item(519,4).NumberFormat = $numFmt
$ws.Cells.Item(520,3).Value = 4
$ws.Cells.Item(520,4).Formula = "=AVERAGE(C514:C520)"
$ws.Cells.Item(520,4).Font.Color = 0
$ws.Cells.Item(520,4).NumberFormat = $numFmt
$ws.Cells.Item(521,3).Value = 4
$ws.Cells.Item(521,4).Formula = "=AVERAGE(C515:C521)"
$ws.Cells.Item(521,4).Font.Color = 0
$ws.Cells.Item(521,4).NumberFormat = $numFmt

# ---------------------------------------------------------------
# 3) Trailing date-only rows 522-541 on sheet4
# ---------------------------------------------------------------
$ws.Cells.Item(522,1).Value = 44420
$ws.Cells.Item(522,1).NumberFormat = $dateFmt
$ws.Cells.Item(523,1).Value = 44421
$ws.Cells.Item(523,1).NumberFormat = $dateFmt
$ws.Cells.Item(524,1).Value = 44422
$ws.Cells.Item(524,1).NumberFormat = $dateFmt
$ws.Cells.Item(525,1).Value = 44423
$ws.Cells.Item(525,1).NumberFormat = $dateFmt
$ws.Cells.Item(526,1).Value = 44424
$ws.Cells.Item(526,1).NumberFormat = $dateFmt
$ws.Cells.Item(527,1).Value = 44425
$ws.Cells.Item(527,1).NumberFormat = $dateFmt
$ws.Cells.Item(528,1).Value = 44426
$ws.Cells.Item(528,1).NumberFormat = $dateFmt
$ws.Cells.Item(529,1).Value = 44427
$ws.Cells.Item(529,1).NumberFormat = $dateFmt
$ws.Cells.Item(530,1).Value = 44428
$ws.Cells.Item(530,1).NumberFormat = $dateFmt
$ws.Cells.Item(531,1).Value = 44429
$ws.Cells.Item(531,1).NumberFormat = $dateFmt
$ws.Cells.Item(532,1).Value = 44430
$ws.Cells.Item(532,1).NumberFormat = $dateFmt
$ws.Cells.Item(533,1).Value = 44431
$ws.Cells.Item(533,1).NumberFormat = $dateFmt
$ws.Cells.Item(534,1).Value = 44432
$ws.Cells.Item(534,1).NumberFormat = $dateFmt
$ws.Cells.Item(535,1).Value = 44433
$ws.Cells.Item(535,1).NumberFormat = $dateFmt
$ws.Cells.Item(536,1).Value = 44434
$ws.Cells.Item(536,1).NumberFormat = $dateFmt
$ws.Cells.Item(537,1).Value = 44435
$ws.Cells.Item(537,1).NumberFormat = $dateFmt
$ws.Cells.Item(538,1).Value = 44436
$ws.Cells.Item(538,1).NumberFormat = $dateFmt
$ws.Cells.Item(539,1).Value = 44437
$ws.Cells.Item(539,1).NumberFormat = $dateFmt
$ws.Cells.Item(540,1).Value = 44438
$ws.Cells.Item(540,1).NumberFormat = $dateFmt
$ws.Cells.Item(541,1).Value = 44439
$ws.Cells.Item(541,1).NumberFormat = $dateFmt

# ---------------------------------------------------------------
# 4) Selection / active-sheet bookkeeping (match saved view state)
#    Sheet3 ("Dimessi   Guariti") is the tab left active/selected.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("A521:D521").Select()
$ws2.Range("A509:D521").Select()
$ws4.Range("A509:D521").Select()
$ws3.Range("A509:D521").Select()

